$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 1.67637130870356

$ws.Range("B3").Value = 0.00009552326474482342
$ws.Range("C3").Value = 0.04103571897497393
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 14.05458458003941
